# [ST 20 April 2018] - Updated code for developers tab
# Insert a new "DeveloperTabData" worksheet right after "Visualize" (i.e. before
# the current 2nd sheet, "DescriptiveStatistics_Data") and populate it with the
# developer/job-history table shown in the target workbook.

$wb = $excel.ActiveWorkbook

# --- 1. Create & position the new sheet ------------------------------------
$ws = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$ws.Name = "DeveloperTabData"

# --- 2. Write the *new* distinct text values first, in the exact sequence
# they were first introduced in the original workbook, so the rebuilt
# shared-string table lines up with the source file (first use wins the next
# available shared-string index).
$ws.Range("C1").Value = "Shell Script"
$ws.Range("D1").Value = "Suyog Talathi"
$ws.Range("B6").Value = "testautocomponent_628539"
$ws.Range("G1").Value = "Created"
$ws.Range("B4").Value = "testautocomponent_166931"
$ws.Range("B5").Value = "testautocomponent_3973701"
$ws.Range("G2").Value = "Updated"
$ws.Range("B3").Value = "testautocomponent_643036"
$ws.Range("B2").Value = "testautocomponent_912424"
$ws.Range("B1").Value = "testautocomponent_961785"

# --- 3. Fill in the remaining values (numbers + repeated strings) ----------
# Row 1
$ws.Range("A1").Value = 320
$ws.Range("E1").Value = 43209.541562500002
$ws.Range("F1").NumberFormat = "m/d/yy h:mm"

# Row 2
$ws.Range("A2").Value = 319
$ws.Range("C2").Value = "Shell Script"
$ws.Range("D2").Value = "Suyog Talathi"
$ws.Range("E2").Value = 43209.521828703706
$ws.Range("F2").Value = 43209.537789351853

# Row 3
$ws.Range("A3").Value = 318
$ws.Range("C3").Value = "Shell Script"
$ws.Range("D3").Value = "Suyog Talathi"
$ws.Range("E3").Value = 43209.490682870368
$ws.Range("F3").NumberFormat = "m/d/yy h:mm"
$ws.Range("G3").Value = "Created"

# Row 4
$ws.Range("C4").Value = "Shell Script"
$ws.Range("D4").Value = "Suyog Talathi"
$ws.Range("A4").Value = 317
$ws.Range("E4").Value = 43209.484513888892
$ws.Range("F4").NumberFormat = "m/d/yy h:mm"
$ws.Range("G4").Value = "Created"

# Row 5
$ws.Range("A5").Value = 316
$ws.Range("C5").Value = "Shell Script"
$ws.Range("D5").Value = "Suyog Talathi"
$ws.Range("E5").Value = 43209.420578703706
$ws.Range("F5").Value = 43209.454317129632
$ws.Range("G5").Value = "Updated"

# Row 6
$ws.Range("A6").Value = 315
$ws.Range("C6").Value = "Shell Script"
$ws.Range("D6").Value = "Suyog Talathi"
$ws.Range("E6").Value = 43209.416608796295
$ws.Range("G6").Value = "Created"

# --- 4. Number formats for the date/time columns ---------------------------
$ws.Range("E1:E6").NumberFormat = "m/d/yy h:mm"
$ws.Range("F1:F5").NumberFormat = "m/d/yy h:mm"

# --- 5. Column widths (best-fit approximation) ------------------------------
$ws.Columns.Item(1).ColumnWidth = 4
$ws.Columns.Item(2).ColumnWidth = 25.25
$ws.Columns.Item(3).ColumnWidth = 9.917
$ws.Columns.Item(4).ColumnWidth = 11.75
$ws.Columns.Item(5).ColumnWidth = 13.584
$ws.Columns.Item(7).ColumnWidth = 8

# --- 6. Final selection (matches the saved view state in the workbook) -----
$ws.Range("D6").Select()
